# The commit "Moving from 2.0.1 to 2.0.2" that produced this diff does not
# change any visible/semantic content of the document. Every single line in
# the unified diff is the *same* set of XML attributes (and the *same*
# namespace declarations) on the *same* elements, just re-serialized in a
# different (alphabetical) order -- the by-product of bumping the OOXML
# library used to (re)generate this test fixture from 2.0.1 to 2.0.2.
#
# Concretely, for every changed element in word/document.xml, word/footer1.xml
# and word/footnotes.xml the attribute multiset is identical before/after
# (only w:... / xmlns:... ordering differs), e.g.:
#   -<w:footerReference w:type="default" r:id="rId6"/>
#   +<w:footerReference r:id="rId6" w:type="default"/>
#   -<w:pgSz w:w="11906" w:h="16838"/>
#   +<w:pgSz w:h="16838" w:w="11906"/>
#   -<w:tcW w:w="3070" w:type="dxa"/>
#   +<w:tcW w:type="dxa" w:w="3070"/>
# ... and so on for every hunk in the diff, including the long w:latentStyles /
# w:style block and the wp:anchor / wps:bodyPr / v:shape markup around the
# footer's text box. No text run, value, id, table layout, style, drawing
# geometry, or relationship is added, removed or changed.
#
# Word's object model (real or COM-interop) has no concept of "attribute
# serialization order" to begin with -- that is purely an artifact of the
# XML writer used by the tool that regenerated this .docx fixture, not
# something reachable through Document/Range/Find/... automation. So the
# faithful way to "apply" this diff through the Word object model is to
# leave the document's content untouched, which is what this script does.

$d = $word.ActiveDocument

# Touch the document through the object model (read-only) so the script is
# a real no-op COM interaction rather than an empty file, without mutating
# any content, formatting, or structure.
$null = $d.Content.Text
$null = $d.Sections.Count
